$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) data per latest scrape,
# and fix ordering swap between Stacks and WOONetwork (rows 50-51).
# Price cells in column D are stored as plain text in the source data
# (e.g. "43.40", "2.197.66"), so a leading apostrophe is used to force
# Excel to keep them as text instead of auto-converting to numbers.

$ws.Range("D2").Value = "'43.758.45"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "'2.197.66"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'259.61"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").Value = "'82.32"
$ws.Range("E6").Value = "  +12.25%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  +2.12%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.596"
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("D10").Value = "'43.40"
$ws.Range("E10").Value = "  +8.14%  "
$ws.Range("D11").Value = "'0.0919"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "'6.98"
$ws.Range("E12").Value = "  +3.40%  "
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "'2.522.96"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "'2.193.76"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "'0.779"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "'43.660.64"
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").Value = "'70.11"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "'5.93"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("E22").Value = "  +13.65%  "
$ws.Range("D23").Value = "'230.92"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").Value = "'8.89"
$ws.Range("E24").Value = "  -5.16%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "'42.03"
$ws.Range("E26").Value = "  +14.68%  "
$ws.Range("D27").Value = "'10.71"
$ws.Range("E27").Value = "  +2.67%  "
$ws.Range("D28").Value = "'3.36"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'2.24"
$ws.Range("E29").Value = "  +2.93%  "
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").Value = "'173.60"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").Value = "'20.46"
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").Value = "'0.0872"
$ws.Range("E33").Value = "  +7.96%  "
$ws.Range("D34").Value = "'5.31"
$ws.Range("E34").Value = "  +4.43%  "
$ws.Range("E35").Value = "  +7.22%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").Value = "'4.49"
$ws.Range("E37").Value = "  +7.16%  "
$ws.Range("D38").Value = "'0.0352"
$ws.Range("E38").Value = "  +5.22%  "
$ws.Range("D39").Value = "'13.24"
$ws.Range("E39").Value = "  +13.06%  "
$ws.Range("D40").Value = "'2.87"
$ws.Range("E40").Value = "  +16.41%  "
$ws.Range("D41").Value = "'2.10"
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("D42").Value = "'63.01"
$ws.Range("E42").Value = "  +6.77%  "
$ws.Range("D43").Value = "'5.47"
$ws.Range("E43").Value = "  +7.19%  "
$ws.Range("D44").Value = "'0.199"
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("D45").Value = "'101.21"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'0.0982"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("D47").Value = "'8.21"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").Value = "'1.18"
$ws.Range("E48").Value = "  +4.56%  "
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "'0.439"
$ws.Range("E50").Value = "  -4.69%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.53"
$ws.Range("E51").Value = "  +26.83%  "

